$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43 (ALC)
$ws.Range("H43").Value = 1954.2307
$ws.Range("I43").Value = 1575
$ws.Range("J43").Value = 2122.7778
$ws.Range("K43").Value = 1575
$ws.Range("L43").Value = 2122.7778
$ws.Range("M43").Value = -1506
$ws.Range("N43").Value = -2260.7778

# Row 69 (ALC)
$ws.Range("H69").Value = 3437.375
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3437.375
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 10312.125
$ws.Range("N69").Value = -12060.125
$ws.Range("M69").ClearContents()

# Row 72 (ALC)
$ws.Range("H72").Value = 3437.375
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3437.375
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 30936.375
$ws.Range("N72").Value = -39672.375
$ws.Range("M72").ClearContents()

# Row 74 (ALC)
$ws.Range("H74").Value = 4961.4443
$ws.Range("I74").Value = 4719
$ws.Range("J74").Value = 5155.4
$ws.Range("K74").Value = 4719
$ws.Range("L74").Value = 5155.4
$ws.Range("M74").Value = -3783
$ws.Range("N74").Value = -7027.4

# Row 77 (ALC)
$ws.Range("H77").Value = 4961.4443
$ws.Range("I77").Value = 4719
$ws.Range("J77").Value = 5155.4
$ws.Range("K77").Value = 23595
$ws.Range("L77").Value = 25777
$ws.Range("M77").Value = -18915
$ws.Range("N77").Value = -35137

# Row 106 (ALC)
$ws.Range("H106").Value = 26090640
$ws.Range("I106").Value = 33337094
$ws.Range("J106").Value = 3400
$ws.Range("K106").Value = 33337094
$ws.Range("L106").Value = 3400
$ws.Range("M106").Value = -33336463
$ws.Range("N106").Value = -4662

# Row 112 (ALC)
$ws.Range("H112").Value = 8959.6875
$ws.Range("J112").Value = 9219.678
$ws.Range("L112").Value = 27659.034
$ws.Range("N112").Value = -29875.034

# Row 116 (ALC)
$ws.Range("H116").Value = 5266.8066
$ws.Range("I116").Value = 7859.706
$ws.Range("J116").Value = 2118.2856
$ws.Range("K116").Value = 7859.706
$ws.Range("L116").Value = 2118.2856
$ws.Range("M116").Value = -4417.706
$ws.Range("N116").Value = -9002.285599999999

# Row 132 (ALC)
$ws.Range("H132").Value = 1988.9634
$ws.Range("I132").Value = 2089.121
$ws.Range("J132").Value = 1575.8125
$ws.Range("K132").Value = 6267.363
$ws.Range("L132").Value = 4727.4375
$ws.Range("M132").Value = -3737.363
$ws.Range("N132").Value = -9787.4375

# Row 138 (ALC)
$ws.Range("H138").Value = 1473.97
$ws.Range("I138").Value = 571.3143
$ws.Range("J138").Value = 1960.0154
$ws.Range("K138").Value = 1713.9429
$ws.Range("L138").Value = 5880.0462
$ws.Range("M138").Value = 3426.0571
$ws.Range("N138").Value = -16160.0462

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 783672.4399999999
$ws.Range("I32").Value = 895191.25
$ws.Range("J32").Value = 23317
$ws.Range("K32").Value = 895191.25
$ws.Range("L32").Value = 23317
$ws.Range("M32").Value = -894904.25
$ws.Range("N32").Value = -23891

# Row 61 (ARM)
$ws.Range("H61").Value = 4567707.5
$ws.Range("I61").Value = 5377580.5
$ws.Range("J61").Value = 2968.7273
$ws.Range("K61").Value = 5377580.5
$ws.Range("L61").Value = 2968.7273
$ws.Range("M61").Value = -5377368.5
$ws.Range("N61").Value = -3392.7273

# Row 102 (ARM)
$ws.Range("H102").Value = 1848.7778
$ws.Range("I102").Value = 1817.375
$ws.Range("J102").Value = 2100
$ws.Range("K102").Value = 1817.375
$ws.Range("L102").Value = 2100
$ws.Range("M102").Value = -195.375
$ws.Range("N102").Value = -5344

# Row 122 (ARM)
$ws.Range("H122").Value = 1597.4286
$ws.Range("I122").Value = 1261.4783
$ws.Range("J122").Value = 3142.8
$ws.Range("K122").Value = 3784.4349
$ws.Range("L122").Value = 9428.400000000001
$ws.Range("M122").Value = -1334.4349
$ws.Range("N122").Value = -14328.4

# Row 132 (ARM)
$ws.Range("H132").Value = 4605.511
$ws.Range("I132").Value = 5271.696
$ws.Range("J132").Value = 3909.0454
$ws.Range("K132").Value = 15815.088
$ws.Range("L132").Value = 11727.1362
$ws.Range("M132").Value = -13285.088
$ws.Range("N132").Value = -16787.1362

# Row 136 (ARM)
$ws.Range("H136").Value = 4567707.5
$ws.Range("I136").Value = 5377580.5
$ws.Range("J136").Value = 2968.7273
$ws.Range("K136").Value = 16132741.5
$ws.Range("L136").Value = 8906.1819
$ws.Range("M136").Value = -16130191.5
$ws.Range("N136").Value = -14006.1819

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (BSM)
$ws.Range("H99").Value = 1046.6666
$ws.Range("I99").Value = 918.1818
$ws.Range("J99").Value = 1400
$ws.Range("K99").Value = 918.1818
$ws.Range("L99").Value = 1400
$ws.Range("M99").Value = 579.8182
$ws.Range("N99").Value = -4396

# Row 105 (BSM)
$ws.Range("H105").Value = 12501431
$ws.Range("I105").Value = 13890279
$ws.Range("K105").Value = 13890279
$ws.Range("M105").Value = -13888532

# Row 134 (BSM)
$ws.Range("H134").Value = 1529.9032
$ws.Range("I134").Value = 1458.9615
$ws.Range("K134").Value = 4376.8845
$ws.Range("M134").Value = -1841.8845

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 3872.8513
$ws.Range("I31").Value = 1251.3903
$ws.Range("J31").Value = 7129.8184
$ws.Range("K31").Value = 1251.3903
$ws.Range("L31").Value = 7129.8184
$ws.Range("M31").Value = -956.3903
$ws.Range("N31").Value = -7719.8184

# Row 34 (CRP)
$ws.Range("H34").Value = 3872.8513
$ws.Range("I34").Value = 1251.3903
$ws.Range("J34").Value = 7129.8184
$ws.Range("K34").Value = 1251.3903
$ws.Range("L34").Value = 7129.8184
$ws.Range("M34").Value = -1049.3903
$ws.Range("N34").Value = -7533.8184

# Row 58 (CRP)
$ws.Range("H58").Value = 793.8706
$ws.Range("I58").Value = 612.8982999999999
$ws.Range("J58").Value = 1204.5385
$ws.Range("K58").Value = 612.8982999999999
$ws.Range("L58").Value = 1204.5385
$ws.Range("M58").Value = -409.8982999999999
$ws.Range("N58").Value = -1610.5385

# Row 122 (CRP)
$ws.Range("H122").Value = 1600.8438
$ws.Range("J122").Value = 1794.55
$ws.Range("L122").Value = 5383.65
$ws.Range("N122").Value = -10283.65

# Row 132 (CRP)
$ws.Range("H132").Value = 4066749.2
$ws.Range("I132").Value = 1501.6451
$ws.Range("K132").Value = 4504.9353
$ws.Range("M132").Value = -1974.9353

# Row 136 (CRP)
$ws.Range("H136").Value = 793.8706
$ws.Range("I136").Value = 612.8982999999999
$ws.Range("J136").Value = 1204.5385
$ws.Range("K136").Value = 1838.6949
$ws.Range("L136").Value = 3613.6155
$ws.Range("M136").Value = 711.3051
$ws.Range("N136").Value = -8713.6155

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 1410.9117
$ws.Range("I5").Value = 491.8
$ws.Range("J5").Value = 2136.5264
$ws.Range("K5").Value = 1475.4
$ws.Range("L5").Value = 6409.5792
$ws.Range("M5").Value = -1363.4
$ws.Range("N5").Value = -6633.5792

# Row 86 (CUL)
$ws.Range("H86").Value = 427.66666
$ws.Range("I86").Value = 480
$ws.Range("J86").Value = 401.5
$ws.Range("K86").Value = 1440
$ws.Range("L86").Value = 1204.5
$ws.Range("M86").Value = -254
$ws.Range("N86").Value = -3576.5

# Row 89 (CUL)
$ws.Range("H89").Value = 427.66666
$ws.Range("I89").Value = 480
$ws.Range("J89").Value = 401.5
$ws.Range("K89").Value = 4320
$ws.Range("L89").Value = 3613.5
$ws.Range("M89").Value = 1608
$ws.Range("N89").Value = -15469.5

# Row 118 (CUL)
$ws.Range("H118").Value = 2708
$ws.Range("I118").Value = 835
$ws.Range("J118").Value = 2996.1538
$ws.Range("K118").Value = 2505
$ws.Range("L118").Value = 8988.4614
$ws.Range("M118").Value = -1262
$ws.Range("N118").Value = -11474.4614

# Row 122 (CUL)
$ws.Range("H122").Value = 2381.2856
$ws.Range("I122").Value = 382.30304
$ws.Range("J122").Value = 5249.391
$ws.Range("K122").Value = 3440.72736
$ws.Range("L122").Value = 47244.519
$ws.Range("M122").Value = -990.7273599999999
$ws.Range("N122").Value = -52144.519

# Row 135 (CUL)
$ws.Range("H135").Value = 1410.9117
$ws.Range("I135").Value = 491.8
$ws.Range("J135").Value = 2136.5264
$ws.Range("K135").Value = 4426.2
$ws.Range("L135").Value = 19228.7376
$ws.Range("M135").Value = -1891.2
$ws.Range("N135").Value = -24298.7376

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (GSM)
$ws.Range("H122").Value = 2900
$ws.Range("I122").Value = 2833.3333
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8499.999899999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6049.999899999999
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
# Row 122 (LTW)
$ws.Range("H122").Value = 3611.9302
$ws.Range("I122").Value = 2865.125
$ws.Range("J122").Value = 4555.263
$ws.Range("K122").Value = 8595.375
$ws.Range("L122").Value = 13665.789
$ws.Range("M122").Value = -6145.375
$ws.Range("N122").Value = -18565.789

# Row 132 (LTW)
$ws.Range("H132").Value = 2127.1843
$ws.Range("I132").Value = 1861.25
$ws.Range("J132").Value = 2871.8
$ws.Range("K132").Value = 5583.75
$ws.Range("L132").Value = 8615.400000000001
$ws.Range("M132").Value = -3053.75
$ws.Range("N132").Value = -13675.4

$ws = $wb.Worksheets.Item("WVR")
# Row 96 (WVR)
$ws.Range("H96").Value = 4772.375
$ws.Range("I96").Value = 3470
$ws.Range("J96").Value = 6074.75
$ws.Range("K96").Value = 3470
$ws.Range("L96").Value = 6074.75
$ws.Range("M96").Value = -2097
$ws.Range("N96").Value = -8820.75

# Row 122 (WVR)
$ws.Range("H122").Value = 3794.7856
$ws.Range("I122").Value = 3243
$ws.Range("K122").Value = 9729
$ws.Range("M122").Value = -7279

# Row 126 (WVR)
$ws.Range("H126").Value = 1390.4375
$ws.Range("I126").Value = 1283.1333
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 3849.3999
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -1379.3999
$ws.Range("N126").Value = -13940

# Row 132 (WVR)
$ws.Range("H132").Value = 5610801
$ws.Range("I132").Value = 1887.4117
$ws.Range("K132").Value = 5662.2351
$ws.Range("M132").Value = -3132.2351

# Row 136 (WVR)
$ws.Range("H136").Value = 2054.24
$ws.Range("I136").Value = 1956.8379
$ws.Range("J136").Value = 2331.4614
$ws.Range("K136").Value = 5870.5137
$ws.Range("L136").Value = 6994.3842
$ws.Range("M136").Value = -3320.5137
$ws.Range("N136").Value = -12094.3842
